$d = $word.ActiveDocument

# Locate the paragraph that ends with the sentence about "uma boa solução."
$rng = $d.Content
$found = $rng.Find.Execute("garantir a convergência a uma boa solução.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the anchor paragraph text."
}

$anchorPara = $rng.Paragraphs(1)

# Insert a brand-new paragraph right after the anchor paragraph. Word
# automatically carries over the paragraph formatting (pPr) from the
# paragraph it was split from, matching the target spacing/indent/jc/sz.
$anchorPara.Range.InsertParagraphAfter()

$newPara = $anchorPara.Next()
$newPara.Range.Text = "A inicialização adequada dos pesos é crucial não apenas em modelos lineares simples, como no exemplo, mas também em redes neurais profundas. Pesos mal inicializados (por exemplo, todos iguais a zero ou com magnitudes desbalanceadas) podem levar a problemas como gradientes vanishing ou exploding, dificultando o treinamento. No contexto de fine-tuning, o modelo parte de um ponto já adaptado a padrões similares, reduzindo o risco de cair em mínimos locais ruins. Assim, a escolha inteligente dos valores iniciais, aliada a taxas de aprendizado adequadas, é essencial para eficiência e robustez em otimização, seja em modelos simples ou em arquiteturas complexas."

Write-Output "Inserted new paragraph. Document now has $($d.Paragraphs.Count) paragraphs."
